# Apply the "xlsx to lua string format" edit:
# In the LanguageSetting sheet, convert the old "{0}"/"{1}" style
# placeholder strings to the new "s%" style placeholder strings.

$wb = $excel.ActiveWorkbook

$langSheet = $wb.Worksheets.Item("LanguageSetting")

# Update the four localized format strings.
$langSheet.Range("B6").Value = "阵营:s%"
$langSheet.Range("B7").Value = "细胞数:s%}/s%"
$langSheet.Range("B8").Value = "已占领:s%/s%"
$langSheet.Range("B9").Value = "分支数:s%/s%"

# Reflect the user's final selection / active sheet state: the
# LanguageSetting tab became the active tab, with B11 selected.
$langSheet.Activate()
$langSheet.Range("B11").Select()
